$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Best Accuracy" header (E1) and add two new headers (F1, G1)
$ws.Range("E1").Value = "Best Accuracy (Val-Split)"
$ws.Range("F1").Value = "Best Accuracy (Val-Excel)"
$ws.Range("G1").Value = "Best Accuracy (Val-Excel) uses test size of 0.00001"

# New value cell F6 containing a pre-formatted text report
$f6 = "accuracy    precision-neg    recall-neg    f1-neg`n----------  ---------------  ------------  --------`n86.95%      85.55%           88.89%        87.19%"
$ws.Range("F6").Value = $f6
$ws.Range("F6").NumberFormat = "0.00%"
$ws.Range("F6").WrapText = $true

# Column widths (values chosen so the stored OOXML width lands on the target)
$ws.Range("E1").ColumnWidth = 20.0
$ws.Range("F1").ColumnWidth = 39.5
$ws.Range("G1").ColumnWidth = 42.666666666666664

# Selection moved to F7 in the saved file
$ws.Range("F7").Select()
